# Applies the "Wijzigingen en aanvullingen prognoses." edit:
#  - Rename the single worksheet from "prog002" to "data"
#  - Add page setup info (A4 paper, portrait orientation) to the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet
$ws.Name = "data"

# Configure page setup: A4 paper size, portrait orientation
$ps = $ws.PageSetup
$ps.PaperSize = 9      # xlPaperA4
$ps.Orientation = 1    # xlPortrait
